$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values
$ws.Range("A2").Value = "test2"
$ws.Range("C2").Value = "safari"

# Add new row 3
$ws.Range("A3").Value = "test2"
$ws.Range("B3").Value = "TD1"
$ws.Range("C3").Value = "firefox"
$ws.Range("D3").Value = "Fail"

# Add new row 4
$ws.Range("A4").Value = "Test1"
$ws.Range("B4").Value = "TD1"
$ws.Range("C4").Value = "firefox"
$ws.Range("D4").Value = "Fail"
